$wb = $excel.ActiveWorkbook

# --- open_tasks sheet: remove the "refactor ffn, dc, m_nb" row (old row 9) ---
$ws1 = $wb.Worksheets.Item("open_tasks")
$ws1.Rows.Item(9).Delete()

# --- done_tasks sheet: update hour totals and append the two new task rows ---
$ws2 = $wb.Worksheets.Item("done_tasks")

# Preference Extraction row
$ws2.Range("B8").Value = 210
$ws2.Range("F8").Value = 180

# Design of FSM row
$ws2.Range("C9").Value = 90
$ws2.Range("D9").Value = 210
$ws2.Range("F9").Value = 30
$ws2.Rows.Item(9).RowHeight = 18.75

# Insert "Implementing FSM" row by copying the formatting of row 9
$ws2.Rows.Item(9).Copy()
$ws2.Rows.Item(10).Insert()
$ws2.Range("A10").Value = "Implementing FSM"
$ws2.Range("B10").ClearContents()
$ws2.Range("C10").ClearContents()
$ws2.Range("D10").ClearContents()
$ws2.Range("E10").Value = 90
$ws2.Range("F10").Value = 30
$ws2.Rows.Item(10).RowHeight = 18.75

# Insert "refactor ffn, dc, m_nb" row (plain right-aligned number format, no border)
$ws2.Rows.Item(11).Insert()
$ws2.Range("E2").Copy()
$ws2.Range("B11:F11").PasteSpecial(-4122)
$ws2.Range("A11").Value = "refactor ffn, dc, m_nb"
$ws2.Range("F11").Value = 30
$ws2.Rows.Item(11).RowHeight = 18.75

# Insert trailing blank row (mirrors the blank row at the bottom of open_tasks);
# explicitly re-paste row 9's formatting since Insert() otherwise copies the row above (11)
$ws2.Rows.Item(12).Insert()
$ws2.Range("A9:F9").Copy()
$ws2.Range("A12:F12").PasteSpecial(-4122)
$ws2.Range("B12:F12").ClearContents()
$ws2.Range("A12").Value = "'"
$ws2.Rows.Item(12).RowHeight = 19.5

$excel.DisplayAlerts = $false
